$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 15 additional rows after row 28 to make room for the expanded player list (rows 29-43)
$ws.Rows("29:43").Insert()

# Clear existing contents for the full data range; styles/formatting on rows 2-28 are preserved,
# new rows 29-43 currently carry the style copied down from row 28 and will be recolored below.
$ws.Range("A2:E28").ClearContents()

# Reassign the alternating green/yellow row-block fill color for every 3-row player block
# (Group1 / Group2 / Difference), since players have been re-sorted alphabetically and
# several new (non-COVID-era) players have been appended.
$ws.Range("A2:E2").Interior.Color = 13434828
$ws.Range("A3:E3").Interior.Color = 13434828
$ws.Range("A4:E4").Interior.Color = 13434828
$ws.Range("A5:E5").Interior.Color = 12451839
$ws.Range("A6:E6").Interior.Color = 12451839
$ws.Range("A7:E7").Interior.Color = 12451839
$ws.Range("A8:E8").Interior.Color = 13434828
$ws.Range("A9:E9").Interior.Color = 13434828
$ws.Range("A10:E10").Interior.Color = 13434828
$ws.Range("A11:E11").Interior.Color = 12451839
$ws.Range("A12:E12").Interior.Color = 12451839
$ws.Range("A13:E13").Interior.Color = 12451839
$ws.Range("A14:E14").Interior.Color = 13434828
$ws.Range("A15:E15").Interior.Color = 13434828
$ws.Range("A16:E16").Interior.Color = 13434828
$ws.Range("A17:E17").Interior.Color = 12451839
$ws.Range("A18:E18").Interior.Color = 12451839
$ws.Range("A19:E19").Interior.Color = 12451839
$ws.Range("A20:E20").Interior.Color = 13434828
$ws.Range("A21:E21").Interior.Color = 13434828
$ws.Range("A22:E22").Interior.Color = 13434828
$ws.Range("A23:E23").Interior.Color = 12451839
$ws.Range("A24:E24").Interior.Color = 12451839
$ws.Range("A25:E25").Interior.Color = 12451839
$ws.Range("A26:E26").Interior.Color = 13434828
$ws.Range("A27:E27").Interior.Color = 13434828
$ws.Range("A28:E28").Interior.Color = 13434828
$ws.Range("A29:E29").Interior.Color = 12451839
$ws.Range("A30:E30").Interior.Color = 12451839
$ws.Range("A31:E31").Interior.Color = 12451839
$ws.Range("A32:E32").Interior.Color = 13434828
$ws.Range("A33:E33").Interior.Color = 13434828
$ws.Range("A34:E34").Interior.Color = 13434828
$ws.Range("A35:E35").Interior.Color = 12451839
$ws.Range("A36:E36").Interior.Color = 12451839
$ws.Range("A37:E37").Interior.Color = 12451839
$ws.Range("A38:E38").Interior.Color = 13434828
$ws.Range("A39:E39").Interior.Color = 13434828
$ws.Range("A40:E40").Interior.Color = 13434828
$ws.Range("A41:E41").Interior.Color = 12451839
$ws.Range("A42:E42").Interior.Color = 12451839
$ws.Range("A43:E43").Interior.Color = 12451839

# Write player name, season-group label, and the three stat columns for every row
$ws.Range("A2").Value = "Austin Hooper"
$ws.Range("B2").Value = "Group1"
$ws.Range("C2").Value = 9.700000000000001
$ws.Range("D2").Value = 6.666666666666667
$ws.Range("E2").Value = 56.33333333333334
$ws.Range("A3").Value = "Austin Hooper"
$ws.Range("B3").Value = "Group2"
$ws.Range("C3").Value = 10.26666666666667
$ws.Range("D3").Value = 7.666666666666667
$ws.Range("E3").Value = 61.33333333333334
$ws.Range("A4").Value = "Austin Hooper"
$ws.Range("B4").Value = "Difference"
$ws.Range("C4").Value = 0.5666666666666664
$ws.Range("D4").Value = 1.0
$ws.Range("E4").Value = 5.0
$ws.Range("A5").Value = "Chris Manhertz"
$ws.Range("B5").Value = "Group1"
$ws.Range("C5").Value = 10.5
$ws.Range("D5").Value = 8.466666666666667
$ws.Range("E5").Value = 72.23333333333333
$ws.Range("A6").Value = "Chris Manhertz"
$ws.Range("B6").Value = "Group2"
$ws.Range("C6").Value = 8.333333333333334
$ws.Range("D6").Value = 5.600000000000001
$ws.Range("E6").Value = 54.16666666666666
$ws.Range("A7").Value = "Chris Manhertz"
$ws.Range("B7").Value = "Difference"
$ws.Range("C7").Value = -2.166666666666666
$ws.Range("D7").Value = -2.866666666666666
$ws.Range("E7").Value = -18.06666666666667
$ws.Range("A8").Value = "Dallas Goedert"
$ws.Range("B8").Value = "Group1"
$ws.Range("C8").Value = 12.23333333333333
$ws.Range("D8").Value = 8.666666666666666
$ws.Range("E8").Value = 61.96666666666667
$ws.Range("A9").Value = "Dallas Goedert"
$ws.Range("B9").Value = "Group2"
$ws.Range("C9").Value = 11.53333333333333
$ws.Range("D9").Value = 8.933333333333332
$ws.Range("E9").Value = 61.86666666666667
$ws.Range("A10").Value = "Dallas Goedert"
$ws.Range("B10").Value = "Difference"
$ws.Range("C10").Value = -0.7000000000000011
$ws.Range("D10").Value = 0.2666666666666657
$ws.Range("E10").Value = -0.09999999999999432
$ws.Range("A11").Value = "Geoff Swaim"
$ws.Range("B11").Value = "Group1"
$ws.Range("C11").Value = 7.0
$ws.Range("D11").Value = 5.333333333333333
$ws.Range("E11").Value = 44.76666666666667
$ws.Range("A12").Value = "Geoff Swaim"
$ws.Range("B12").Value = "Group2"
$ws.Range("C12").Value = 7.066666666666666
$ws.Range("D12").Value = 6.366666666666667
$ws.Range("E12").Value = 70.06666666666666
$ws.Range("A13").Value = "Geoff Swaim"
$ws.Range("B13").Value = "Difference"
$ws.Range("C13").Value = 0.06666666666666643
$ws.Range("D13").Value = 1.033333333333334
$ws.Range("E13").Value = 25.29999999999999
$ws.Range("A14").Value = "Gerald Everett"
$ws.Range("B14").Value = "Group1"
$ws.Range("C14").Value = 10.4
$ws.Range("D14").Value = 7.033333333333334
$ws.Range("E14").Value = 52.29999999999999
$ws.Range("A15").Value = "Gerald Everett"
$ws.Range("B15").Value = "Group2"
$ws.Range("C15").Value = 7.399999999999999
$ws.Range("D15").Value = 5.033333333333334
$ws.Range("E15").Value = 40.46666666666667
$ws.Range("A16").Value = "Gerald Everett"
$ws.Range("B16").Value = "Difference"
$ws.Range("C16").Value = -3.000000000000001
$ws.Range("D16").Value = -2.0
$ws.Range("E16").Value = -11.83333333333332
$ws.Range("A17").Value = "Hayden Hurst"
$ws.Range("B17").Value = "Group1"
$ws.Range("C17").Value = 10.1
$ws.Range("D17").Value = 7.5
$ws.Range("E17").Value = 59.46666666666667
$ws.Range("A18").Value = "Hayden Hurst"
$ws.Range("B18").Value = "Group2"
$ws.Range("C18").Value = 9.1
$ws.Range("D18").Value = 5.833333333333333
$ws.Range("E18").Value = 49.06666666666666
$ws.Range("A19").Value = "Hayden Hurst"
$ws.Range("B19").Value = "Difference"
$ws.Range("C19").Value = -1.0
$ws.Range("D19").Value = -1.666666666666667
$ws.Range("E19").Value = -10.40000000000001
$ws.Range("A20").Value = "John Mundt"
$ws.Range("B20").Value = "Group1"
$ws.Range("C20").Value = 9.6
$ws.Range("D20").Value = 8.666666666666666
$ws.Range("E20").Value = 80.96666666666667
$ws.Range("A21").Value = "John Mundt"
$ws.Range("B21").Value = "Group2"
$ws.Range("C21").Value = 8.333333333333334
$ws.Range("D21").Value = 6.5
$ws.Range("E21").Value = 54.43333333333334
$ws.Range("A22").Value = "John Mundt"
$ws.Range("B22").Value = "Difference"
$ws.Range("C22").Value = -1.266666666666666
$ws.Range("D22").Value = -2.166666666666666
$ws.Range("E22").Value = -26.53333333333333
$ws.Range("A23").Value = "Marcedes Lewis"
$ws.Range("B23").Value = "Group1"
$ws.Range("C23").Value = 10.13333333333333
$ws.Range("D23").Value = 7.366666666666667
$ws.Range("E23").Value = 57.93333333333333
$ws.Range("A24").Value = "Marcedes Lewis"
$ws.Range("B24").Value = "Group2"
$ws.Range("C24").Value = 6.766666666666667
$ws.Range("D24").Value = 5.399999999999999
$ws.Range("E24").Value = 67.13333333333334
$ws.Range("A25").Value = "Marcedes Lewis"
$ws.Range("B25").Value = "Difference"
$ws.Range("C25").Value = -3.366666666666668
$ws.Range("D25").Value = -1.966666666666668
$ws.Range("E25").Value = 9.20000000000001
$ws.Range("A26").Value = "Mike Gesicki"
$ws.Range("B26").Value = "Group1"
$ws.Range("C26").Value = 11.73333333333333
$ws.Range("D26").Value = 7.233333333333334
$ws.Range("E26").Value = 52.03333333333333
$ws.Range("A27").Value = "Mike Gesicki"
$ws.Range("B27").Value = "Group2"
$ws.Range("C27").Value = 9.966666666666667
$ws.Range("D27").Value = 6.8
$ws.Range("E27").Value = 51.73333333333333
$ws.Range("A28").Value = "Mike Gesicki"
$ws.Range("B28").Value = "Difference"
$ws.Range("C28").Value = -1.766666666666667
$ws.Range("D28").Value = -0.4333333333333345
$ws.Range("E28").Value = -0.3000000000000043
$ws.Range("A29").Value = "Nick Vannett"
$ws.Range("B29").Value = "Group1"
$ws.Range("C29").Value = 10.43333333333333
$ws.Range("D29").Value = 6.977777777777779
$ws.Range("E29").Value = 52.9
$ws.Range("A30").Value = "Nick Vannett"
$ws.Range("B30").Value = "Group2"
$ws.Range("C30").Value = 6.544444444444444
$ws.Range("D30").Value = 5.333333333333333
$ws.Range("E30").Value = 80.3888888888889
$ws.Range("A31").Value = "Nick Vannett"
$ws.Range("B31").Value = "Difference"
$ws.Range("C31").Value = -3.888888888888889
$ws.Range("D31").Value = -1.644444444444446
$ws.Range("E31").Value = 27.4888888888889
$ws.Range("A32").Value = "Noah Fant"
$ws.Range("B32").Value = "Group1"
$ws.Range("C32").Value = 11.63333333333333
$ws.Range("D32").Value = 7.7
$ws.Range("E32").Value = 47.26666666666667
$ws.Range("A33").Value = "Noah Fant"
$ws.Range("B33").Value = "Group2"
$ws.Range("C33").Value = 11.0
$ws.Range("D33").Value = 8.366666666666667
$ws.Range("E33").Value = 57.43333333333334
$ws.Range("A34").Value = "Noah Fant"
$ws.Range("B34").Value = "Difference"
$ws.Range("C34").Value = -0.6333333333333329
$ws.Range("D34").Value = 0.666666666666667
$ws.Range("E34").Value = 10.16666666666666
$ws.Range("A35").Value = "Travis Kelce"
$ws.Range("B35").Value = "Group1"
$ws.Range("C35").Value = 12.8
$ws.Range("D35").Value = 9.066666666666668
$ws.Range("E35").Value = 63.53333333333333
$ws.Range("A36").Value = "Travis Kelce"
$ws.Range("B36").Value = "Group2"
$ws.Range("C36").Value = 10.43333333333333
$ws.Range("D36").Value = 7.699999999999999
$ws.Range("E36").Value = 61.8
$ws.Range("A37").Value = "Travis Kelce"
$ws.Range("B37").Value = "Difference"
$ws.Range("C37").Value = -2.366666666666667
$ws.Range("D37").Value = -1.366666666666669
$ws.Range("E37").Value = -1.733333333333327
$ws.Range("A38").Value = "Will Dissly"
$ws.Range("B38").Value = "Group1"
$ws.Range("C38").Value = 10.96666666666667
$ws.Range("D38").Value = 9.1
$ws.Range("E38").Value = 63.43333333333334
$ws.Range("A39").Value = "Will Dissly"
$ws.Range("B39").Value = "Group2"
$ws.Range("C39").Value = 10.0
$ws.Range("D39").Value = 8.166666666666666
$ws.Range("E39").Value = 62.5
$ws.Range("A40").Value = "Will Dissly"
$ws.Range("B40").Value = "Difference"
$ws.Range("C40").Value = -0.9666666666666668
$ws.Range("D40").Value = -0.9333333333333336
$ws.Range("E40").Value = -0.9333333333333371
$ws.Range("A41").Value = "Zach Ertz"
$ws.Range("B41").Value = "Group1"
$ws.Range("C41").Value = 10.02222222222222
$ws.Range("D41").Value = 6.055555555555556
$ws.Range("E41").Value = 48.07777777777778
$ws.Range("A42").Value = "Zach Ertz"
$ws.Range("B42").Value = "Group2"
$ws.Range("C42").Value = 8.466666666666667
$ws.Range("D42").Value = 5.8
$ws.Range("E42").Value = 53.56666666666666
$ws.Range("A43").Value = "Zach Ertz"
$ws.Range("B43").Value = "Difference"
$ws.Range("C43").Value = -1.555555555555557
$ws.Range("D43").Value = -0.2555555555555564
$ws.Range("E43").Value = 5.488888888888887
